# "updated activity till excel form"
# Swap the runs (C) / balls (D) figures recorded for Varun Aaron's two
# innings against Rajasthan Royals. Values are stored as text (the sheet
# keeps numberStoredAsText), so force text entry with a leading apostrophe
# to avoid Excel re-typing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'0"
$ws.Range("D2").Value = "'8"

$ws.Range("C3").Value = "'1"
$ws.Range("D3").Value = "'2"
